# Backend form-submission sheet: add a "select domain" column (C) with a
# default "General" choice plus a "Clear selection" pass, and log every
# submission's timestamp in the new column D (shifted out from the old C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: resize/reposition the workbook window to match the saved view
# (mirrors the new xWindow/yWindow/windowWidth/windowHeight in bookViews).
$win = $excel.ActiveWindow
$win.Top = -120
$win.Width = 28110
$win.Height = 16440

# Start from a clean used range so stale cells (old B/C layout) don't linger.
$ws.UsedRange.Clear()

# Row 1 - "sgfdgasd" / "dfgfd" submitted with the default domain ("General")
$ws.Range("A1").Value = "sgfdgasd"
$ws.Range("B1").Value = "dfgfd"
$ws.Range("C1").Value = "General"
$ws.Range("D1").Value = 45187.69875416667
$ws.Range("D1").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 2 - "hello " / "hii" with default domain ("General")
$ws.Range("A2").Value = "hello "
$ws.Range("B2").Value = "hii"
$ws.Range("C2").Value = "General"
$ws.Range("D2").Value = 45187.69890496528
$ws.Range("D2").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 3 - "hello " / "hii" with domain "Admin"
$ws.Range("A3").Value = "hello "
$ws.Range("B3").Value = "hii"
$ws.Range("C3").Value = "Admin"
$ws.Range("D3").Value = 45187.69913512732
$ws.Range("D3").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 4 - "hello " / "hii" with domain "Finance"
$ws.Range("A4").Value = "hello "
$ws.Range("B4").Value = "hii"
$ws.Range("C4").Value = "Finance"
$ws.Range("D4").Value = 45187.69916291667
$ws.Range("D4").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 5 - "hello " / "hii" with domain "News"
$ws.Range("A5").Value = "hello "
$ws.Range("B5").Value = "hii"
$ws.Range("C5").Value = "News"
$ws.Range("D5").Value = 45187.69919403935
$ws.Range("D5").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 6 - "hello " / "hii" with domain "Legal"
$ws.Range("A6").Value = "hello "
$ws.Range("B6").Value = "hii"
$ws.Range("C6").Value = "Legal"
$ws.Range("D6").Value = 45187.69921873264
$ws.Range("D6").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Clear-selection feature: leave the grid focused on the freshly-added block
# (mirrors picking "Select all" then clearing down to a fresh A1 anchor).
$ws.Range("A1:D11").Select() | Out-Null
